# edit.ps1 - applies commit diff "Process actual source document: 何大成理事長採訪稿" to before.docx
# via Word COM-interop object model.
#
# NOTE: this runtime's PowerShell-like engine does not reliably bind named
# (-Param value) arguments to function param() blocks, so every helper below
# is called with positional arguments only.

$d = $word.ActiveDocument

function Replace-ParaText($OldText, $NewText) {
    $found = $d.Content.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $found) {
        throw "Find/Replace failed for text: $OldText"
    }
}

function Find-ParaIndexByPrefix($Prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.StartsWith($Prefix)) {
            return $i
        }
    }
    throw "Paragraph not found with prefix: $Prefix"
}

function Insert-ParasAfterIndex($AfterIndex, $Items) {
    $idx = $AfterIndex
    foreach ($item in $Items) {
        $d.Paragraphs($idx).Range.InsertParagraphAfter() | Out-Null
        $idx = $idx + 1
        $d.Paragraphs($idx).Range.Text = $item
    }
}

# ---------------------------------------------------------------------------
# Step 1: Text replacements for "受訪者經歷" bullets 1-3, and the four article
# sections FirstParagraph / BodyText paragraphs that were rewritten wholesale.
# ---------------------------------------------------------------------------
Replace-ParaText '曾任新北市政府消防局 化災搶救隊員' '財團法人工業技術研究院 正工程師'
Replace-ParaText '參與多起重大化災事故應變' '財團法人工業技術研究院 綠能所 室經理'
Replace-ParaText '推動臺灣化學應變體系建置' '財團法人工業技術研究院 能環所 室主任'
Replace-ParaText '雖然臺灣應變體系起步較晚，但在政府與民間共同努力下，已建立起具有特色的化學應變機制。中華民國化學應變協會理事長何大成，從消防員到應變專家，見證並參與了臺灣化災應變的成長歷程，他相信只要持續培育人才、落實企業責任，臺灣的應變能力將能與國際接軌，共同守護未來。' '猶記得，當初設計出ERIC LOGO時，請西服店老闆幫我繡上，老闆回了我一句話：「少年，我們要縫自己的LOGO，不要老是用國外別人的。」其實那時候我很自豪，因為我可以大聲說：「這是我自己設計的。」雖然目前臺灣的應變體系及訓練制度都是參考國外，但經過多年來政府及民間的努力，儼然已完成一套專屬臺灣的應變體系。一個品牌的LOGO，把要傳達給大家的精神及意涵都濃縮在裡面。這LOGO是我自創的，這也正是我們想成立屬於臺灣自己的『化學應變協會』的初衷。'
Replace-ParaText '臺灣早期在毒化災領域相對落後，何大成理事長回憶起當年出國取經的經歷，深深感受到國際間的差距。「當時我們到美國、德國學習先進的應變技術，發現他們已經有非常完整的體系和設備，而臺灣還在起步階段。」然而，他也認識到不能完全照搬國外經驗，必須因地制宜，發展出適合臺灣的應變模式。' '「提著一卡皮箱走天下」這可以說是何大成理事長一生致力臺灣毒化災應變的寫照。他回顧起最初接觸毒化災領域那段時間，其實蠻孤獨的。起初，臺灣在這塊領域較為落後，必須去跟國外交流學習，但～何大成說：「我們沒有什麼可以給予人家，別人又怎麼願意教予我們呢？」就這樣獨身前往美國及德國這兩個工業大國，跟他們一起學習，以行動來證明自己的價值。'
Replace-ParaText '經過多年努力，臺灣逐漸建立起自己的應變體系，結合國際經驗與在地需求。何理事長表示：「我們學習美國的 ERIC（緊急應變資訊中心）概念，但也發展出符合臺灣產業特性的做法。」這種取鏡世界、惠及臺灣的精神，成為臺灣化學應變發展的重要基石。' '何大成提到，他其實很感謝美國。美國在工業上，已有一兩百年的歷史演變到現在的組織架構。而他們也很慷慨，樂於跟我們分享資訊及設備，更樂於與你討論交流。但是在德國就完全不同，何大成提到曾有一次，他想向德國的朋友索取設備清單，但被婉拒。對方提到：「我可以給你參考，但是～沒有用！因為那是我的。」這樣簡單的一句話震撼何大成！其實對方說得沒有錯。這些設備、這些制度，是適合德國當地，但～拿回臺灣，不同的地理、環境，真的適合嗎？'
Replace-ParaText '在何理事長看來，應變工作的核心永遠是人命安全。「以人命為基準，這是不能妥協的原則。」他強調，企業在追求營運效率的同時，更應該善盡社會責任，建立完善的應變機制。' '何大成提到，臺灣法規及訓練制度已臻完善，但對於第三方應變單位的支援服務，仍在起步階段，許多業者還無法認同「使用者付費」這個觀念，會認為這是政府單位的責任。何大成舉例在國外，如：美國杜邦公司。他們每項產品都會主動擬定SDS(化學品安全技術說明書)，甚至自主性建立聯防架構。他們不覺得應變單是政府單位的工作，也是企業社會責任（CSR）。所以美國是業者先帶起，接著換政府制定相關法規。在法規制定完成後，民間應變體系也早就跟上。'
Replace-ParaText '「企業責任才是應變的保護傘。」何理事長指出，許多化學災害的發生，往往與企業的風險管理不足有關。因此，CERA 積極推動企業建立自主應變能力，不僅是為了符合法規要求，更是為了保護員工、社區民眾與環境。唯有企業真正認知到應變的重要性，並投入資源建置，才能在災害發生時有效因應。' '反觀臺灣，法規及相關辦法已走出一條屬於自己的路，可是國際上如美國CHEMTREC及德國TUIS等知名第三方應變單位，難道臺灣不需要嗎？「安全文化」是臺灣企業尚缺的概念。也就因這樣的契機，何大成在民國104年與有志一同的業界夥伴成立「中華民國化學應變協會CERA」。在去年化學署頒布「毒性及關注化學物質環境事故專業應變諮詢機關（構）認證及管理辦法」後，CERA也是第一家受到國家認證的應變機關（構）。'
Replace-ParaText '「有能力就有責任，不能不做。」這是何理事長的信念，也是他投身應變工作數十年的動力。從消防隊員時期的第一線搶救，到現在協會理事長的角色，他始終堅持在人命為優先的前提下，盡自己所能協助建立更安全的環境。' '從ERIC到CERA那充滿設計性的LOGO是何大成理事長不假他人之手親手繪製，就如同這LOGO的創作，他知道這體制不得不建立，有「使命感」才能將「安全、技術與夥伴」結合在一起。'
Replace-ParaText '何理事長分享，臺灣的應變人員往往面臨資源有限、人力不足的挑戰，但大家仍然堅守崗位。「每次出勤，我都告訴自己，一定要把弟兄安全帶回來。」這份責任感，不僅是對同仁的承諾，也是對社會的承諾。他期許更多有能力的專業人士能投入應變領域，共同守護臺灣的安全。' '何大成談到，最開始的草創，由工業技術研究院ERIC、強本、輝宇、仲信、柏旭等公司的加入，現在已可達到全台北、中、南三區兩小時內到達應變現場，現階段也正與東部地區業者洽談中，相信不久CERA的應變涵蓋網就可觸及全台，目前CERA與國際同步提供會員業者從Level1專業諮詢、應變諮詢提供等，Level2教育訓練服務及事故處理紀錄等，到Level3的事故廠ERC緊急應變中心及派遣專業應變人員到事故現場等，就是為了實現國內一系列完整的專業應變服務。'
Replace-ParaText '經過多年耕耘，臺灣的應變體系已經從萌芽期進入成長期。何理事長欣慰地表示：「應變的種子已經種下，也開始成長，但需要大家共同來澆灌扶持。」' '「你！要是沒有使命感，就不要進這個門。」做應變是一個志業，因為有志同道合的朋友，這條路走起來才會覺得快樂也才能走得長遠。這是何大成對應變的態度。CERA作為臺灣第一個取得政府認證的第三方應變單位只是一個種子，從國外開始逐步借鏡到現在，CERA已經屬於臺灣自有的”品種”。何大成希望，這只是一個開始，希望未來有更多的第三方應變單位加入，所謂的應變，是要資源共享；應變，也不是商業上的競爭，更不是經濟上的衝突，未來甚至可以合作聯盟，國內有更多的應變機關（構），據點越多越能加快應變的時效性，安全的應變才是大家共通的目標。'
Replace-ParaText 'CERA 致力於培育應變人才，透過訓練課程、經驗交流、國際合作等方式，讓更多人具備應變能力。何理事長也呼籲政府、企業、學術界攜手合作，共同建構更完善的應變網絡。「臺灣雖小，但我們的應變能力不能小。」他相信，只要持續投入、不斷精進，臺灣的應變品牌將能在國際上佔有一席之地，與世界各國一同守護未來。' '何大成提到，國外地廣人稀，應變手段與臺灣截然不同。如VR應用在應變這一塊，我們與美國、德國、日本等國家就像是在競爭賽跑。臺灣的優勢就是科技知識，因為臺灣訓練的土地不夠，藉由這些科技，反而讓我們沉浸在完全不同以往的訓練方式，不管是場地大小的選定、毒化物數值的調整，都可以無限大不再有侷限，在ERIC的努力下，我們走在這領域的尖端，不再同於以往這是屬於臺灣的應變。'

# ---------------------------------------------------------------------------
# Step 2: Paragraph insertions (new bullets / new body paragraphs).
# Each anchor is located by its (already-replaced) text so that indices do not
# need to be hard-coded and stay correct even as earlier insertions shift them.
# ---------------------------------------------------------------------------
$anchorIdx = Find-ParaIndexByPrefix '財團法人工業技術研究院 能環所 室主任'
Insert-ParasAfterIndex $anchorIdx @('財團法人工業技術研究院 環安中心 室經理', '財團法人工業技術研究院 緊急應變諮詢中心（ERIC）資深諮詢員', '環保署環境事故應變諮詢中心 資深諮詢員', '環保署毒災應變諮詢中心 資深諮詢員', '中華民國化學應變協會（CERA）理事', '仁德醫護管理專科學校職安系兼任講師', '明新科技大學化學工程系兼任講師', '國科會科學紀錄資料中心約聘人員', '逢甲大學化學工程系專任助教')

$anchorIdx = Find-ParaIndexByPrefix '何大成提到，他其實很感謝美國。美國在工業上，已有一兩百年的歷史演變到現在的組織架構。而他們也很慷慨，樂於跟我們分享資訊及設備，更樂於與你討論交流。但是在德國就完全不同，何大成提到曾有一次，他想向德國的朋友索取設備清單，但被婉拒。對方提到：「我可以給你參考，但是～沒有用！因為那是我的。」這樣簡單的一句話震撼何大成！其實對方說得沒有錯。這些設備、這些制度，是適合德國當地，但～拿回臺灣，不同的地理、環境，真的適合嗎？'
Insert-ParasAfterIndex $anchorIdx @('何大成更回想到，在美國受訓時教官提到：所有的應變訓練課綱，都要因地制宜，不可能一套方式適用全部。就因為這些經歷，這才促進了下一階段。何大成自傲說，過去臺灣的應變體制及訓練，都是源自美國NFPA及HazMat，但今日如果跟國外比較，我們已經不輸給其他國家。')

$anchorIdx = Find-ParaIndexByPrefix '何大成談到，最開始的草創，由工業技術研究院ERIC、強本、輝宇、仲信、柏旭等公司的加入，現在已可達到全台北、中、南三區兩小時內到達應變現場，現階段也正與東部地區業者洽談中，相信不久CERA的應變涵蓋網就可觸及全台，目前CERA與國際同步提供會員業者從Level1專業諮詢、應變諮詢提供等，Level2教育訓練服務及事故處理紀錄等，到Level3的事故廠ERC緊急應變中心及派遣專業應變人員到事故現場等，就是為了實現國內一系列完整的專業應變服務。'
Insert-ParasAfterIndex $anchorIdx @('而CERA為了能提升設備的整合應用率，更將應變器材如：SCAB、防護衣等，品牌規格統一化，應變時會員成員間器材可交互支援，在訓練上也更能達到事半功倍的效果，說到訓練，何大成更提到，要在應變現場能最快速有效作為就是平時的扎實訓練，每年CERA會召集會員至少集訓一次以上，而各會員每一季也會自主性訓練，如應變技能、法規及偵測儀器操作等，為的就是在平時累積應變量能。而為了能與政府部門相互搭配，CERA除了參與大型的實兵及沙盤演練外，聯防無預警測試也不缺席，在無預警的時間壓力測試下才能實證應變的量能。', '何大成更特別提到在ERIC及CERA相互配合下，已經是一個完美的整合，從前端的諮詢、資訊到專家級的指揮官派遣及整個應變團隊，甚至到虛擬實境及擴增實境的訓練等等都能獨立完成，有國外專家曾經說過：ERIC在全世界是個非常特殊的單位，從諮詢、訓練到研發及應變，沒有類似一個單位可以做到這樣，也當然沒有其他單位可以介入的角色。如今ERIC及CERA在臺灣法規的驗證下，各司諮詢及應變的兩個獨立角色，並且互相協作配合，將讓臺灣應變體系有更大的發展空間。', '回頭談到在應變的這一塊，何大成在外是鐵漢，但在內最不捨的卻是讓母親擔心受怕。他回想起之前台中一起工廠事故，火勢非常大。回到家時，恰巧母親正在看著這起事故新聞，抬頭看見何大成，以緊張的神情說：「這地方很危險，不要去。」但～何大成不敢說實話，只能內心默默說道：「已經去過現場了。」在家裡，何大成是老母親唯一的兒子；在外面，何大成是應變機關（構）的靈魂人物，因為他有一個使命感不得不做。就是「要把在外應變的人員，一個一個安全的帶回來，這才是一個完整的應變。」', '「自助」、「互助」與「聯防」是應變的三道防線，第三方應變機構就是事故廠面對自身應變能量不足時，提供專業技術人員及設備到場支援的應變單位，何大成舉例：假如我有六個據點，你有八個據點而他有四個據點，這樣全部加起來就有十八個，越多的待命應變支援點，在災故時才能快速抵達並穩定現場，達到有效的減災作用，而CERA就是本著這個概念建立起來，以最快速最專業的能力，協助會員共同解決問題，讓毒化災應變這責任不是由公部門，也不是由業者，更不是CERA單獨承擔，而是三方共同努力達到三贏的局面。')

$anchorIdx = Find-ParaIndexByPrefix '何大成提到，國外地廣人稀，應變手段與臺灣截然不同。如VR應用在應變這一塊，我們與美國、德國、日本等國家就像是在競爭賽跑。臺灣的優勢就是科技知識，因為臺灣訓練的土地不夠，藉由這些科技，反而讓我們沉浸在完全不同以往的訓練方式，不管是場地大小的選定、毒化物數值的調整，都可以無限大不再有侷限，在ERIC的努力下，我們走在這領域的尖端，不再同於以往這是屬於臺灣的應變。'
Insert-ParasAfterIndex $anchorIdx @('臺灣在政府部門的法規健全下，相關應變體系也逐漸完整，CERA僅只是一個起始點，現在臺灣的應變儼然發展成獨立的一支。何大成表示：「只要有心，在政府單位與民間業界的共同努力下，相信整個體系完善是值得期待的。」希望能有更多相同理念的企業一同加入，讓第三方應變單位在臺灣遍地開花，撐起臺灣應變的保護傘，守護臺灣也鼎足世界，共同守護地球。')

Write-Host "Done. Final paragraph count: $($d.Paragraphs.Count)"
